$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 934, shifting the existing rows 934-967 down to 935-968
$ws.Rows.Item(934).Insert()

# Populate the newly inserted row 934 with the new weekly record
$ws.Cells.Item(934, 1).Value = 6
$ws.Cells.Item(934, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(934, 3).Value = "Metropolitana"
$ws.Cells.Item(934, 4).Value = 44939
$ws.Cells.Item(934, 5).Value = 13
$ws.Cells.Item(934, 6).Value = 100112003
$ws.Cells.Item(934, 7).Value = "Ajo"
$ws.Cells.Item(934, 8).Value = "Chino"
$ws.Cells.Item(934, 9).Value = "Primera"
$ws.Cells.Item(934, 10).Value = 1100
$ws.Cells.Item(934, 11).Value = 12000
$ws.Cells.Item(934, 12).Value = 13000
$ws.Cells.Item(934, 13).Value = 12455
$ws.Cells.Item(934, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(934, 15).Value = "China"
$ws.Cells.Item(934, 16).Value = 1246
$ws.Cells.Item(934, 17).Value = 10
$ws.Cells.Item(934, 18).Value = "Hortaliza"
